$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424226571106146
$ws.Range("D2").Value = 0.02057710033479765
$ws.Range("E2").Value = 0.1500079334831597
$ws.Range("F2").Value = 0.6474669960146215
$ws.Range("G2").Value = 0.002410213022355554
$ws.Range("K2").Value = 0.504451890541219
$ws.Range("M2").Value = 0.2594698726378297
$ws.Range("N2").Value = 1.420477081238527
$ws.Range("O2").Value = 2.172609407967229
$ws.Range("B3").Value = 0.1329384921194645
$ws.Range("D3").Value = 0.01909246799107223
$ws.Range("E3").Value = 0.1417913037927434
$ws.Range("F3").Value = 0.6410921036915411
$ws.Range("G3").Value = 0.002412831595020937
$ws.Range("K3").Value = 0.4399034382151399
$ws.Range("M3").Value = 0.2310841556912564
$ws.Range("N3").Value = 1.437070428696972
$ws.Range("O3").Value = 2.16501735059515
$ws.Range("B4").Value = 0.1271846491612649
$ws.Range("D4").Value = 0.0181732226077429
$ws.Range("E4").Value = 0.1368639139511032
$ws.Range("F4").Value = 0.637616863124606
$ws.Range("G4").Value = 0.002414525806350379
$ws.Range("K4").Value = 0.4001375809899912
$ws.Range("M4").Value = 0.2137145984690889
$ws.Range("N4").Value = 1.447799345253058
$ws.Range("O4").Value = 2.161794295690783
$ws.Range("B5").Value = 0.1248575453243177
$ws.Range("D5").Value = 0.01779671434477592
$ws.Range("E5").Value = 0.1348853227037594
$ws.Range("F5").Value = 0.6363109712809418
$ws.Range("G5").Value = 0.002415238002186849
$ws.Range("K5").Value = 0.3838999374272021
$ws.Range("M5").Value = 0.2066512616597294
$ws.Range("N5").Value = 1.452307372438872
$ws.Range("O5").Value = 2.160842219427224
$ws.Range("B6").Value = 0.1244722013834405
$ws.Range("D6").Value = 0.01773408081660222
$ws.Range("E6").Value = 0.1345585458441647
$ws.Range("F6").Value = 0.6361007881133816
$ws.Range("G6").Value = 0.002415357579604339
$ws.Range("K6").Value = 0.3812017318617222
$ws.Range("M6").Value = 0.2054792988486938
$ws.Range("N6").Value = 1.45306413494972
$ws.Range("O6").Value = 2.160705941231953
$ws.Range("B7").Value = 0.1271531934079349
$ws.Range("D7").Value = 0.01816815258331417
$ws.Range("E7").Value = 0.1368371113965736
$ws.Range("F7").Value = 0.6375988049370491
$ws.Range("G7").Value = 0.002414535322981083
$ws.Range("K7").Value = 0.3999187260968142
$ws.Range("M7").Value = 0.2136192795584648
$ws.Range("N7").Value = 1.447859591823924
$ws.Range("O7").Value = 2.161779993091557
$ws.Range("B8").Value = 0.139138193429531
$ws.Range("D8").Value = 0.02006680464755561
$ws.Range("E8").Value = 0.1471502823457698
$ws.Range("F8").Value = 0.6451777697305303
$ws.Range("G8").Value = 0.002411098014067097
$ws.Range("K8").Value = 0.4822235399192607
$ws.Range("M8").Value = 0.2496701217400243
$ws.Range("N8").Value = 1.426086271726062
$ws.Range("O8").Value = 2.169692869728863
$ws.Range("B9").Value = 0.1631860419953313
$ws.Range("D9").Value = 0.02372836676921963
$ws.Range("E9").Value = 0.1683192111666898
$ws.Range("F9").Value = 0.6635287709480622
$ws.Range("G9").Value = 0.002405039984598325
$ws.Range("K9").Value = 0.6425470573905727
$ws.Range("M9").Value = 0.3208421747779653
$ws.Range("N9").Value = 1.387677852510821
$ws.Range("O9").Value = 2.196645093461228
$ws.Range("B10").Value = 0.1811803485757792
$ws.Range("D10").Value = 0.02638008553910254
$ws.Range("E10").Value = 0.1844653212161305
$ws.Range("F10").Value = 0.6791484555532605
$ws.Range("G10").Value = 0.002401001031881985
$ws.Range("K10").Value = 0.7596626024628108
$ws.Range("M10").Value = 0.3734353091154006
$ws.Range("N10").Value = 1.362073174822497
$ws.Range("O10").Value = 2.223453544629876
$ws.Range("B11").Value = 0.1894361076589917
$ws.Range("D11").Value = 0.02757791570711277
$ws.Range("E11").Value = 0.1919432563060823
$ws.Range("F11").Value = 0.6867207250180201
$ws.Range("G11").Value = 0.002399252149669942
$ws.Range("K11").Value = 0.8127919408509854
$ws.Range("M11").Value = 0.3974300996903821
$ws.Range("N11").Value = 1.350992421439715
$ws.Range("O11").Value = 2.237178868670583
$ws.Range("B12").Value = 0.1925722799148843
$ws.Range("D12").Value = 0.0280302686957512
$ws.Range("E12").Value = 0.1947943518120141
$ws.Range("F12").Value = 0.6896554112272355
$ws.Range("G12").Value = 0.002398602546081413
$ws.Range("K12").Value = 0.8328889569566797
$ws.Range("M12").Value = 0.4065264780484483
$ws.Range("N12").Value = 1.346877983134483
$ws.Range("O12").Value = 2.242596840798569
$ws.Range("B13").Value = 0.191896411856959
$ws.Range("D13").Value = 0.02793290188909481
$ws.Range("E13").Value = 0.1941794523450469
$ws.Range("F13").Value = 0.6890203821758973
$ws.Range("G13").Value = 0.002398741887696281
$ws.Range("K13").Value = 0.8285616915073604
$ws.Range("M13").Value = 0.404566961961919
$ws.Range("N13").Value = 1.34776047048965
$ws.Range("O13").Value = 2.241420169672779
$ws.Range("B14").Value = 0.1896939251443115
$ws.Range("D14").Value = 0.02761515603486231
$ws.Range("E14").Value = 0.1921774281932187
$ws.Range("F14").Value = 0.6869608153748601
$ws.Range("G14").Value = 0.002399198452956552
$ws.Range("K14").Value = 0.8144457779284551
$ws.Range("M14").Value = 0.3981782622479244
$ws.Range("N14").Value = 1.350652288237256
$ws.Range("O14").Value = 2.237620186841212
$ws.Range("B15").Value = 0.1883461218020841
$ws.Range("D15").Value = 0.02742036541477688
$ws.Range("E15").Value = 0.190953659462906
$ws.Range("F15").Value = 0.6857080300183185
$ws.Range("G15").Value = 0.002399479759433872
$ws.Range("K15").Value = 0.8057964973289131
$ws.Range("M15").Value = 0.3942663114254827
$ws.Range("N15").Value = 1.352434239496617
$ws.Range("O15").Value = 2.23532131597446
$ws.Range("B16").Value = 0.180642209782377
$ws.Range("D16").Value = 0.02630163241595795
$ws.Range("E16").Value = 0.183979317166461
$ws.Range("F16").Value = 0.6786629922658847
$ws.Range("G16").Value = 0.002401117100249322
$ws.Range("K16").Value = 0.7561874502769683
$ws.Range("M16").Value = 0.3718686015113875
$ws.Range("N16").Value = 1.362808742951126
$ws.Range("O16").Value = 2.222587390974297
$ws.Range("B17").Value = 0.1759339258701829
$ws.Range("D17").Value = 0.02561314447977736
$ws.Range("E17").Value = 0.1797350164340727
$ws.Range("F17").Value = 0.674460724179454
$ws.Range("G17").Value = 0.002402144169084951
$ws.Range("K17").Value = 0.72571568926071
$ws.Range("M17").Value = 0.3581462510880655
$ws.Range("N17").Value = 1.369318426135266
$ws.Range("O17").Value = 2.215167754022332
$ws.Range("B18").Value = 0.1732324504798157
$ws.Range("D18").Value = 0.02521635051903104
$ws.Range("E18").Value = 0.1773063100580359
$ws.Range("F18").Value = 0.6720876214495064
$ws.Range("G18").Value = 0.002402743241440729
$ws.Range("K18").Value = 0.7081753367175452
$ws.Range("M18").Value = 0.3502600972068848
$ws.Range("N18").Value = 1.373115982277847
$ws.Range("O18").Value = 2.211044159966434
$ws.Range("B19").Value = 0.1723189171863453
$ws.Range("D19").Value = 0.02508186728476858
$ws.Range("E19").Value = 0.1764861319238449
$ws.Range("F19").Value = 0.6712916715064381
$ws.Range("G19").Value = 0.002402947509398596
$ws.Range("K19").Value = 0.7022341273839174
$ws.Range("M19").Value = 0.3475911103978007
$ws.Range("N19").Value = 1.37441093163377
$ws.Range("O19").Value = 2.209672696971126
$ws.Range("B20").Value = 0.1764344484342217
$ws.Range("D20").Value = 0.02568651753814777
$ws.Range("E20").Value = 0.1801855337263021
$ws.Range("F20").Value = 0.6749035153448659
$ws.Range("G20").Value = 0.002402033974090301
$ws.Range("K20").Value = 0.7289608941373729
$ws.Range("M20").Value = 0.3596063376646441
$ws.Range("N20").Value = 1.368619936288432
$ws.Range("O20").Value = 2.215942681706821
$ws.Range("B21").Value = 0.1903405817485151
$ws.Range("D21").Value = 0.0277085195214184
$ws.Range("E21").Value = 0.192764943773291
$ws.Range("F21").Value = 0.6875639345570619
$ws.Range("G21").Value = 0.002399064005712689
$ws.Range("K21").Value = 0.8185925646142209
$ws.Range("M21").Value = 0.4000545047187956
$ws.Range("N21").Value = 1.349800675969398
$ws.Range("O21").Value = 2.238730346106081
$ws.Range("B22").Value = 0.199486638753001
$ws.Range("D22").Value = 0.02902277991262281
$ws.Range("E22").Value = 0.2010992940663172
$ws.Range("F22").Value = 0.6962301798583468
$ws.Range("G22").Value = 0.002397196724444374
$ws.Range("K22").Value = 0.877043976479456
$ws.Range("M22").Value = 0.4265484357126752
$ws.Range("N22").Value = 1.337976881224243
$ws.Range("O22").Value = 2.254908805037843
$ws.Range("B23").Value = 0.1946000052834478
$ws.Range("D23").Value = 0.02832200437936194
$ws.Range("E23").Value = 0.196640677461005
$ws.Range("F23").Value = 0.6915689447868942
$ws.Range("G23").Value = 0.002398186597828187
$ws.Range("K23").Value = 0.8458593257980169
$ws.Range("M23").Value = 0.412402747406972
$ws.Range("N23").Value = 1.344243912085746
$ws.Range("O23").Value = 2.24615628606486
$ws.Range("B24").Value = 0.176208145412474
$ws.Range("D24").Value = 0.02565334860749857
$ws.Range("E24").Value = 0.1799818193733245
$ws.Range("F24").Value = 0.6747031960308192
$ws.Range("G24").Value = 0.002402083766291737
$ws.Range("K24").Value = 0.7274938046491286
$ws.Range("M24").Value = 0.3589462231588882
$ws.Range("N24").Value = 1.368935552064606
$ws.Range("O24").Value = 2.215591894439655
$ws.Range("B25").Value = 0.1566227105279552
$ws.Range("D25").Value = 0.02274450670763883
$ws.Range("E25").Value = 0.1624895169464082
$ws.Range("F25").Value = 0.6581897739990552
$ws.Range("G25").Value = 0.002406606209473726
$ws.Range("K25").Value = 0.5992923747097905
$ws.Range("M25").Value = 0.3015359413704672
$ws.Range("N25").Value = 1.39760926594735
$ws.Range("O25").Value = 2.188125926967984
